# Apply the textual updates to the "two-digit number divided by one-digit
# number" practice table. The table has 20 rows, but only every 4th row
# (1, 5, 9, 13, 17 in 1-based terms) actually holds the 5 exercise cells;
# the intervening rows are blank spacer rows. We replace the content of
# each of those 25 cells, in document order, with its new value.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$dataRows = @(1, 5, 9, 13, 17)

$newValues = @(
    "11÷7=", "37÷2=", "20÷8=", "17÷7=", "12÷5=",
    "26÷3=", "95÷3=", "23÷7=", "97÷3=", "45÷4=",
    "67÷3=", "81÷7=", "77÷4=", "78÷6=", "49÷2=",
    "80÷6=", "94÷6=", "98÷2=", "53÷3=", "78÷9=",
    "77÷4=", "68÷2=", "33÷4=", "94÷5=", "14÷2="
)

$idx = 0
foreach ($r in $dataRows) {
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($r, $c)
        $rng = $cell.Range
        # Trim the trailing cell-mark / paragraph-mark characters so we
        # only touch the visible text of the cell.
        $rng.MoveEnd(1, -1) | Out-Null
        $rng.Text = $newValues[$idx]
        $idx = $idx + 1
    }
}
